$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"0.053267"
$ws.Range("H2").Value = [double]"0.159801"
$ws.Range("I2").Value = [double]"0.0002370783698475801"
$ws.Range("J2").Value = [double]"0.00023707836984758"
$ws.Range("M2").Value = [double]"16.087096"
$ws.Range("N2").Value = [double]"48.261288"
$ws.Range("O2").Value = [double]"0.1263055268415452"
$ws.Range("P2").Value = [double]"0.1263055268415452"
$ws.Range("Q2").Value = [double]"0.856911342632"
$ws.Range("R2").Value = [double]"7.712202083688"
$ws.Range("S2").Value = [double]"2.994430840633331E-05"
$ws.Range("T2").Value = [double]"2.994430840633331E-05"

$ws.Range("G3").Value = [double]"0.053267"
$ws.Range("H3").Value = [double]"0.159801"
$ws.Range("I3").Value = [double]"0.0002370783698475801"
$ws.Range("J3").Value = [double]"0.00023707836984758"
$ws.Range("O3").Value = [double]"0.7490048915888087"
$ws.Range("P3").Value = [double]"0.7490048915888088"
$ws.Range("Q3").Value = [double]"5.08157325605"
$ws.Range("R3").Value = [double]"45.73415930445"
$ws.Range("S3").Value = [double]"0.0001775728587057382"
$ws.Range("T3").Value = [double]"0.0001775728587057382"

$ws.Range("G4").Value = [double]"0.053267"
$ws.Range("H4").Value = [double]"0.159801"
$ws.Range("I4").Value = [double]"0.0002370783698475801"
$ws.Range("J4").Value = [double]"0.00023707836984758"
$ws.Range("M4").Value = [double]"0.5200936666666667"
$ws.Range("N4").Value = [double]"1.560281"
$ws.Range("O4").Value = [double]"0.004083440825819921"
$ws.Range("P4").Value = [double]"0.004083440825819921"
$ws.Range("Q4").Value = [double]"0.02770382934233333"
$ws.Range("R4").Value = [double]"0.249334464081"
$ws.Range("S4").Value = [double]"9.680954943544431E-07"
$ws.Range("T4").Value = [double]"9.680954943544429E-07"

$ws.Range("G5").Value = [double]"0.053267"
$ws.Range("H5").Value = [double]"0.159801"
$ws.Range("I5").Value = [double]"0.0002370783698475801"
$ws.Range("J5").Value = [double]"0.00023707836984758"
$ws.Range("M5").Value = [double]"14.15205133333333"
$ws.Range("N5").Value = [double]"42.456154"
$ws.Range("O5").Value = [double]"0.1111128011883101"
$ws.Range("P5").Value = [double]"0.1111128011883101"
$ws.Range("Q5").Value = [double]"0.7538373183726667"
$ws.Range("R5").Value = [double]"6.784535865354"
$ws.Range("S5").Value = [double]"2.634244177492282E-05"
$ws.Range("T5").Value = [double]"2.634244177492282E-05"

$ws.Range("G6").Value = [double]"0.053267"
$ws.Range("H6").Value = [double]"0.159801"
$ws.Range("I6").Value = [double]"0.0002370783698475801"
$ws.Range("J6").Value = [double]"0.00023707836984758"
$ws.Range("M6").Value = [double]"0.794831"
$ws.Range("N6").Value = [double]"2.384493"
$ws.Range("O6").Value = [double]"0.006240501592393819"
$ws.Range("P6").Value = [double]"0.006240501592393819"
$ws.Range("Q6").Value = [double]"0.042338262877"
$ws.Range("R6").Value = [double]"0.381044365893"
$ws.Range("S6").Value = [double]"1.479487944555954E-06"
$ws.Range("T6").Value = [double]"1.479487944555954E-06"

$ws.Range("G7").Value = [double]"0.053267"
$ws.Range("H7").Value = [double]"0.159801"
$ws.Range("I7").Value = [double]"0.0002370783698475801"
$ws.Range("J7").Value = [double]"0.00023707836984758"
$ws.Range("M7").Value = [double]"0.4143026666666667"
$ws.Range("N7").Value = [double]"1.242908"
$ws.Range("O7").Value = [double]"0.003252837963122146"
$ws.Range("P7").Value = [double]"0.003252837963122146"
$ws.Range("Q7").Value = [double]"0.02206866014533333"
$ws.Range("R7").Value = [double]"0.198617941308"
$ws.Range("S7").Value = [double]"7.711775216753213E-07"
$ws.Range("T7").Value = [double]"7.711775216753212E-07"

$ws.Range("G8").Value = [double]"9.236317"
$ws.Range("H8").Value = [double]"27.708951"
$ws.Range("I8").Value = [double]"0.0411085846350553"
$ws.Range("J8").Value = [double]"0.0411085846350553"
$ws.Range("M8").Value = [double]"16.087096"
$ws.Range("N8").Value = [double]"48.261288"
$ws.Range("O8").Value = [double]"0.1263055268415452"
$ws.Range("P8").Value = [double]"0.1263055268415452"
$ws.Range("Q8").Value = [double]"148.585518265432"
$ws.Range("R8").Value = [double]"1337.269664388888"
$ws.Range("S8").Value = [double]"0.005192241440040911"
$ws.Range("T8").Value = [double]"0.005192241440040911"

$ws.Range("G9").Value = [double]"9.236317"
$ws.Range("H9").Value = [double]"27.708951"
$ws.Range("I9").Value = [double]"0.0411085846350553"
$ws.Range("J9").Value = [double]"0.0411085846350553"
$ws.Range("O9").Value = [double]"0.7490048915888087"
$ws.Range("P9").Value = [double]"0.7490048915888088"
$ws.Range("Q9").Value = [double]"881.1275546135499"
$ws.Range("R9").Value = [double]"7930.147991521951"
$ws.Range("S9").Value = [double]"0.03079053097794896"
$ws.Range("T9").Value = [double]"0.03079053097794897"

$ws.Range("G10").Value = [double]"9.236317"
$ws.Range("H10").Value = [double]"27.708951"
$ws.Range("I10").Value = [double]"0.0411085846350553"
$ws.Range("J10").Value = [double]"0.0411085846350553"
$ws.Range("M10").Value = [double]"0.5200936666666667"
$ws.Range("N10").Value = [double]"1.560281"
$ws.Range("O10").Value = [double]"0.004083440825819921"
$ws.Range("P10").Value = [double]"0.004083440825819921"
$ws.Range("Q10").Value = [double]"4.803749975025666"
$ws.Range("R10").Value = [double]"43.233749775231"
$ws.Range("S10").Value = [double]"0.0001678644727904583"
$ws.Range("T10").Value = [double]"0.0001678644727904583"

$ws.Range("G11").Value = [double]"9.236317"
$ws.Range("H11").Value = [double]"27.708951"
$ws.Range("I11").Value = [double]"0.0411085846350553"
$ws.Range("J11").Value = [double]"0.0411085846350553"
$ws.Range("M11").Value = [double]"14.15205133333333"
$ws.Range("N11").Value = [double]"42.456154"
$ws.Range("O11").Value = [double]"0.1111128011883101"
$ws.Range("P11").Value = [double]"0.1111128011883101"
$ws.Range("Q11").Value = [double]"130.7128323149393"
$ws.Range("R11").Value = [double]"1176.415490834454"
$ws.Range("S11").Value = [double]"0.00456768999168772"
$ws.Range("T11").Value = [double]"0.00456768999168772"

$ws.Range("G12").Value = [double]"9.236317"
$ws.Range("H12").Value = [double]"27.708951"
$ws.Range("I12").Value = [double]"0.0411085846350553"
$ws.Range("J12").Value = [double]"0.0411085846350553"
$ws.Range("M12").Value = [double]"0.794831"
$ws.Range("N12").Value = [double]"2.384493"
$ws.Range("O12").Value = [double]"0.006240501592393819"
$ws.Range("P12").Value = [double]"0.006240501592393819"
$ws.Range("Q12").Value = [double]"7.341311077426999"
$ws.Range("R12").Value = [double]"66.071799696843"
$ws.Range("S12").Value = [double]"0.0002565381878761187"
$ws.Range("T12").Value = [double]"0.0002565381878761187"

$ws.Range("G13").Value = [double]"9.236317"
$ws.Range("H13").Value = [double]"27.708951"
$ws.Range("I13").Value = [double]"0.0411085846350553"
$ws.Range("J13").Value = [double]"0.0411085846350553"
$ws.Range("M13").Value = [double]"0.4143026666666667"
$ws.Range("N13").Value = [double]"1.242908"
$ws.Range("O13").Value = [double]"0.003252837963122146"
$ws.Range("P13").Value = [double]"0.003252837963122146"
$ws.Range("Q13").Value = [double]"3.826630763278666"
$ws.Range("R13").Value = [double]"34.439676869508"
$ws.Range("S13").Value = [double]"0.0001337195647111277"
$ws.Range("T13").Value = [double]"0.0001337195647111277"

$ws.Range("G14").Value = [double]"212.661977"
$ws.Range("H14").Value = [double]"637.9859310000001"
$ws.Range("I14").Value = [double]"0.9465063704680865"
$ws.Range("J14").Value = [double]"0.9465063704680865"
$ws.Range("M14").Value = [double]"16.087096"
$ws.Range("N14").Value = [double]"48.261288"
$ws.Range("O14").Value = [double]"0.1263055268415452"
$ws.Range("P14").Value = [double]"0.1263055268415452"
$ws.Range("Q14").Value = [double]"3421.113639548792"
$ws.Range("R14").Value = [double]"30790.02275593913"
$ws.Range("S14").Value = [double]"0.1195489857808505"
$ws.Range("T14").Value = [double]"0.1195489857808505"

$ws.Range("G15").Value = [double]"212.661977"
$ws.Range("H15").Value = [double]"637.9859310000001"
$ws.Range("I15").Value = [double]"0.9465063704680865"
$ws.Range("J15").Value = [double]"0.9465063704680865"
$ws.Range("O15").Value = [double]"0.7490048915888087"
$ws.Range("P15").Value = [double]"0.7490048915888088"
$ws.Range("Q15").Value = [double]"20287.55918114255"
$ws.Range("R15").Value = [double]"182588.032630283"
$ws.Range("S15").Value = [double]"0.7089379014005659"
$ws.Range("T15").Value = [double]"0.708937901400566"

$ws.Range("G16").Value = [double]"212.661977"
$ws.Range("H16").Value = [double]"637.9859310000001"
$ws.Range("I16").Value = [double]"0.9465063704680865"
$ws.Range("J16").Value = [double]"0.9465063704680865"
$ws.Range("M16").Value = [double]"0.5200936666666667"
$ws.Range("N16").Value = [double]"1.560281"
$ws.Range("O16").Value = [double]"0.004083440825819921"
$ws.Range("P16").Value = [double]"0.004083440825819921"
$ws.Range("Q16").Value = [double]"110.6041473785123"
$ws.Range("R16").Value = [double]"995.4373264066111"
$ws.Range("S16").Value = [double]"0.003865002755068019"
$ws.Range("T16").Value = [double]"0.003865002755068019"

$ws.Range("G17").Value = [double]"212.661977"
$ws.Range("H17").Value = [double]"637.9859310000001"
$ws.Range("I17").Value = [double]"0.9465063704680865"
$ws.Range("J17").Value = [double]"0.9465063704680865"
$ws.Range("M17").Value = [double]"14.15205133333333"
$ws.Range("N17").Value = [double]"42.456154"
$ws.Range("O17").Value = [double]"0.1111128011883101"
$ws.Range("P17").Value = [double]"0.1111128011883101"
$ws.Range("Q17").Value = [double]"3009.603215152153"
$ws.Range("R17").Value = [double]"27086.42893636938"
$ws.Range("S17").Value = [double]"0.1051689741652895"
$ws.Range("T17").Value = [double]"0.1051689741652895"

$ws.Range("G18").Value = [double]"212.661977"
$ws.Range("H18").Value = [double]"637.9859310000001"
$ws.Range("I18").Value = [double]"0.9465063704680865"
$ws.Range("J18").Value = [double]"0.9465063704680865"
$ws.Range("M18").Value = [double]"0.794831"
$ws.Range("N18").Value = [double]"2.384493"
$ws.Range("O18").Value = [double]"0.006240501592393819"
$ws.Range("P18").Value = [double]"0.006240501592393819"
$ws.Range("Q18").Value = [double]"169.030331840887"
$ws.Range("R18").Value = [double]"1521.272986567983"
$ws.Range("S18").Value = [double]"0.005906674512116987"
$ws.Range("T18").Value = [double]"0.005906674512116987"

$ws.Range("G19").Value = [double]"212.661977"
$ws.Range("H19").Value = [double]"637.9859310000001"
$ws.Range("I19").Value = [double]"0.9465063704680865"
$ws.Range("J19").Value = [double]"0.9465063704680865"
$ws.Range("M19").Value = [double]"0.4143026666666667"
$ws.Range("N19").Value = [double]"1.242908"
$ws.Range("O19").Value = [double]"0.003252837963122146"
$ws.Range("P19").Value = [double]"0.003252837963122146"
$ws.Range("Q19").Value = [double]"88.10642416970533"
$ws.Range("R19").Value = [double]"792.957817527348"
$ws.Range("S19").Value = [double]"0.003078831854195546"
$ws.Range("T19").Value = [double]"0.003078831854195546"

$ws.Range("G20").Value = [double]"2.729417"
$ws.Range("H20").Value = [double]"8.188250999999999"
$ws.Range("I20").Value = [double]"0.01214796652701058"
$ws.Range("J20").Value = [double]"0.01214796652701058"
$ws.Range("M20").Value = [double]"16.087096"
$ws.Range("N20").Value = [double]"48.261288"
$ws.Range("O20").Value = [double]"0.1263055268415452"
$ws.Range("P20").Value = [double]"0.1263055268415452"
$ws.Range("Q20").Value = [double]"43.90839330303199"
$ws.Range("R20").Value = [double]"395.175539727288"
$ws.Range("S20").Value = [double]"0.001534355312247527"
$ws.Range("T20").Value = [double]"0.001534355312247527"

$ws.Range("G21").Value = [double]"2.729417"
$ws.Range("H21").Value = [double]"8.188250999999999"
$ws.Range("I21").Value = [double]"0.01214796652701058"
$ws.Range("J21").Value = [double]"0.01214796652701058"
$ws.Range("O21").Value = [double]"0.7490048915888087"
$ws.Range("P21").Value = [double]"0.7490048915888088"
$ws.Range("Q21").Value = [double]"260.38133237855"
$ws.Range("R21").Value = [double]"2343.43199140695"
$ws.Range("S21").Value = [double]"0.009098886351588033"
$ws.Range("T21").Value = [double]"0.009098886351588033"

$ws.Range("G22").Value = [double]"2.729417"
$ws.Range("H22").Value = [double]"8.188250999999999"
$ws.Range("I22").Value = [double]"0.01214796652701058"
$ws.Range("J22").Value = [double]"0.01214796652701058"
$ws.Range("M22").Value = [double]"0.5200936666666667"
$ws.Range("N22").Value = [double]"1.560281"
$ws.Range("O22").Value = [double]"0.004083440825819921"
$ws.Range("P22").Value = [double]"0.004083440825819921"
$ws.Range("Q22").Value = [double]"1.419552495392333"
$ws.Range("R22").Value = [double]"12.775972458531"
$ws.Range("S22").Value = [double]"4.960550246708882E-05"
$ws.Range("T22").Value = [double]"4.960550246708882E-05"

$ws.Range("G23").Value = [double]"2.729417"
$ws.Range("H23").Value = [double]"8.188250999999999"
$ws.Range("I23").Value = [double]"0.01214796652701058"
$ws.Range("J23").Value = [double]"0.01214796652701058"
$ws.Range("M23").Value = [double]"14.15205133333333"
$ws.Range("N23").Value = [double]"42.456154"
$ws.Range("O23").Value = [double]"0.1111128011883101"
$ws.Range("P23").Value = [double]"0.1111128011883101"
$ws.Range("Q23").Value = [double]"38.62684949407266"
$ws.Range("R23").Value = [double]"347.641645446654"
$ws.Range("S23").Value = [double]"0.001349794589557972"
$ws.Range("T23").Value = [double]"0.001349794589557972"

$ws.Range("G24").Value = [double]"2.729417"
$ws.Range("H24").Value = [double]"8.188250999999999"
$ws.Range("I24").Value = [double]"0.01214796652701058"
$ws.Range("J24").Value = [double]"0.01214796652701058"
$ws.Range("M24").Value = [double]"0.794831"
$ws.Range("N24").Value = [double]"2.384493"
$ws.Range("O24").Value = [double]"0.006240501592393819"
$ws.Range("P24").Value = [double]"0.006240501592393819"
$ws.Range("Q24").Value = [double]"2.169425243527"
$ws.Range("R24").Value = [double]"19.524827191743"
$ws.Range("S24").Value = [double]"7.58094044561563E-05"
$ws.Range("T24").Value = [double]"7.58094044561563E-05"

$ws.Range("G25").Value = [double]"2.729417"
$ws.Range("H25").Value = [double]"8.188250999999999"
$ws.Range("I25").Value = [double]"0.01214796652701058"
$ws.Range("J25").Value = [double]"0.01214796652701058"
$ws.Range("M25").Value = [double]"0.4143026666666667"
$ws.Range("N25").Value = [double]"1.242908"
$ws.Range("O25").Value = [double]"0.003252837963122146"
$ws.Range("P25").Value = [double]"0.003252837963122146"
$ws.Range("Q25").Value = [double]"1.130804741545333"
$ws.Range("R25").Value = [double]"10.177242673908"
$ws.Range("S25").Value = [double]"3.95153666937971E-05"
$ws.Range("T25").Value = [double]"3.95153666937971E-05"
